$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "52.247.75"
$ws.Range("E2").Value = "  +2.02%  "
$ws.Range("D3").Value = "2.796.14"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "346.62"
$ws.Range("E5").Value = "  +3.97%  "
Set-TextValue $ws.Range("D6") "116.42"
$ws.Range("E6").Value = "  +1.14%  "
Set-TextValue $ws.Range("D7") "0.551"
$ws.Range("E7").Value = "  +3.85%  "
$ws.Range("E8").Value = "  -0.07%  "
Set-TextValue $ws.Range("D9") "0.590"
$ws.Range("E9").Value = "  +3.17%  "
Set-TextValue $ws.Range("D10") "43.09"
$ws.Range("E10").Value = "  +4.14%  "
Set-TextValue $ws.Range("D11") "0.0857"
$ws.Range("E11").Value = "  +3.50%  "
Set-TextValue $ws.Range("D12") "20.15"
$ws.Range("E12").Value = "  -0.02%  "
$ws.Range("E13").Value = "  +1.67%  "
Set-TextValue $ws.Range("D14") "7.86"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").Value = "3.238.51"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").Value = "2.779.08"
$ws.Range("E16").Value = "  +1.66%  "
Set-TextValue $ws.Range("D17") "0.894"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "52.156.28"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("E19").Value = "  +7.74%  "
$ws.Range("E20").Value = "  +4.17%  "
$ws.Range("E21").Value = "  -1.21%  "
$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +2.12%  "
Set-TextValue $ws.Range("D23") "70.17"
$ws.Range("E23").Value = "  -0.01%  "
Set-TextValue $ws.Range("D24") "270.09"
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("E25").Value = "  +4.99%  "
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  -0.15%  "
Set-TextValue $ws.Range("D28") "10.22"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  +0.64%  "
$ws.Range("E30").Value = "  -0.15%  "
Set-TextValue $ws.Range("D31") "35.15"
$ws.Range("E31").Value = "  -1.34%  "
$ws.Range("E32").Value = "  +0.50%  "
Set-TextValue $ws.Range("D33") "5.71"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D34") "0.0421"
$ws.Range("E34").Value = "  +19.48%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.0827"
$ws.Range("E35").Value = "  +0.17%  "
Set-TextValue $ws.Range("D36") "2.12"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("E37").Value = "  +0.09%  "
Set-TextValue $ws.Range("D38") "18.92"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("E40").Value = "  +0.08%  "
Set-TextValue $ws.Range("D41") "2.70"
$ws.Range("E41").Value = "  +21.10%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "23.67"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D43") "127.79"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("E44").Value = "  +1.65%  "
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "2.073.19"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  +2.29%  "
Set-TextValue $ws.Range("D49") "0.974"
$ws.Range("E49").Value = "  +14.87%  "
Set-TextValue $ws.Range("D50") "5.52"
$ws.Range("E50").Value = "  -0.36%  "
Set-TextValue $ws.Range("D51") "8.96"
$ws.Range("E51").Value = "  -1.12%  "
